$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: model name headers (including new Neural Network - MPL column)
$ws.Cells.Item(1, 2).Value = 'K-NN'
$ws.Cells.Item(1, 4).Value = 'K-NN Centroid'
$ws.Cells.Item(1, 6).Value = 'Decision Tree'
$ws.Cells.Item(1, 8).Value = 'SVM-SVC'
$ws.Cells.Item(1, 10).Value = 'Naive Bayes'
$ws.Cells.Item(1, 12).Value = 'Random Forest'
$ws.Cells.Item(1, 14).Value = 'Neural Network - MPL'

# Row 2: Average / Standard Deviation sub-headers
$ws.Cells.Item(2, 1).Value = 'Users'
$ws.Cells.Item(2, 2).Value = 'Average'
$ws.Cells.Item(2, 3).Value = 'Standard Deviation'
$ws.Cells.Item(2, 4).Value = 'Average'
$ws.Cells.Item(2, 5).Value = 'Standard Deviation'
$ws.Cells.Item(2, 6).Value = 'Average'
$ws.Cells.Item(2, 7).Value = 'Standard Deviation'
$ws.Cells.Item(2, 8).Value = 'Average'
$ws.Cells.Item(2, 9).Value = 'Standard Deviation'
$ws.Cells.Item(2, 10).Value = 'Average'
$ws.Cells.Item(2, 11).Value = 'Standard Deviation'
$ws.Cells.Item(2, 12).Value = 'Average'
$ws.Cells.Item(2, 13).Value = 'Standard Deviation'
$ws.Cells.Item(2, 14).Value = 'Average'
$ws.Cells.Item(2, 15).Value = 'Standard Deviation'

# Row 3: data (percentages as text, leading spaces preserved)
$ws.Cells.Item(3, 2).Value = '   80.00%'
$ws.Cells.Item(3, 3).Value = '   17.92%'
$ws.Cells.Item(3, 4).Value = '   55.38%'
$ws.Cells.Item(3, 5).Value = '   23.83%'
$ws.Cells.Item(3, 6).Value = '   72.31%'
$ws.Cells.Item(3, 7).Value = '   16.94%'
$ws.Cells.Item(3, 8).Value = '   80.00%'
$ws.Cells.Item(3, 9).Value = '    4.65%'
$ws.Cells.Item(3, 10).Value = '   81.54%'
$ws.Cells.Item(3, 11).Value = '   18.52%'
$ws.Cells.Item(3, 12).Value = '   83.08%'
$ws.Cells.Item(3, 13).Value = '   11.13%'
$ws.Cells.Item(3, 14).Value = '   80.00%'
$ws.Cells.Item(3, 15).Value = '    4.65%'

# Row 4: data (percentages as text, leading spaces preserved)
$ws.Cells.Item(4, 2).Value = '   87.69%'
$ws.Cells.Item(4, 3).Value = '    9.67%'
$ws.Cells.Item(4, 4).Value = '   69.23%'
$ws.Cells.Item(4, 5).Value = '   14.44%'
$ws.Cells.Item(4, 6).Value = '   89.23%'
$ws.Cells.Item(4, 7).Value = '   16.52%'
$ws.Cells.Item(4, 8).Value = '   80.00%'
$ws.Cells.Item(4, 9).Value = '    4.65%'
$ws.Cells.Item(4, 10).Value = '   87.69%'
$ws.Cells.Item(4, 11).Value = '   16.78%'
$ws.Cells.Item(4, 12).Value = '   86.15%'
$ws.Cells.Item(4, 13).Value = '   15.63%'
$ws.Cells.Item(4, 14).Value = '   80.00%'
$ws.Cells.Item(4, 15).Value = '    4.65%'

# Row 5: data (percentages as text, leading spaces preserved)
$ws.Cells.Item(5, 2).Value = '   86.15%'
$ws.Cells.Item(5, 3).Value = '   19.73%'
$ws.Cells.Item(5, 4).Value = '   78.46%'
$ws.Cells.Item(5, 5).Value = '   27.61%'
$ws.Cells.Item(5, 6).Value = '   89.23%'
$ws.Cells.Item(5, 7).Value = '   15.86%'
$ws.Cells.Item(5, 8).Value = '   80.00%'
$ws.Cells.Item(5, 9).Value = '    4.65%'
$ws.Cells.Item(5, 10).Value = '   81.54%'
$ws.Cells.Item(5, 11).Value = '   14.33%'
$ws.Cells.Item(5, 12).Value = '   86.15%'
$ws.Cells.Item(5, 13).Value = '   19.57%'
$ws.Cells.Item(5, 14).Value = '   80.00%'
$ws.Cells.Item(5, 15).Value = '    4.65%'

# Row 6: data (percentages as text, leading spaces preserved)
$ws.Cells.Item(6, 2).Value = '   90.77%'
$ws.Cells.Item(6, 3).Value = '    9.79%'
$ws.Cells.Item(6, 4).Value = '   81.54%'
$ws.Cells.Item(6, 5).Value = '   16.86%'
$ws.Cells.Item(6, 6).Value = '   80.00%'
$ws.Cells.Item(6, 7).Value = '   12.66%'
$ws.Cells.Item(6, 8).Value = '   80.00%'
$ws.Cells.Item(6, 9).Value = '    4.65%'
$ws.Cells.Item(6, 10).Value = '   78.46%'
$ws.Cells.Item(6, 11).Value = '   19.35%'
$ws.Cells.Item(6, 12).Value = '   86.15%'
$ws.Cells.Item(6, 13).Value = '   10.55%'
$ws.Cells.Item(6, 14).Value = '   84.62%'
$ws.Cells.Item(6, 15).Value = '   13.99%'

# Row 7: data (percentages as text, leading spaces preserved)
$ws.Cells.Item(7, 2).Value = '   80.00%'
$ws.Cells.Item(7, 3).Value = '   15.10%'
$ws.Cells.Item(7, 4).Value = '   63.08%'
$ws.Cells.Item(7, 5).Value = '   22.58%'
$ws.Cells.Item(7, 6).Value = '   70.77%'
$ws.Cells.Item(7, 7).Value = '   19.00%'
$ws.Cells.Item(7, 8).Value = '   80.00%'
$ws.Cells.Item(7, 9).Value = '    4.65%'
$ws.Cells.Item(7, 10).Value = '   84.62%'
$ws.Cells.Item(7, 11).Value = '   13.99%'
$ws.Cells.Item(7, 12).Value = '   78.46%'
$ws.Cells.Item(7, 13).Value = '    9.90%'
$ws.Cells.Item(7, 14).Value = '   80.00%'
$ws.Cells.Item(7, 15).Value = '    4.65%'
